$wb = $excel.ActiveWorkbook

# ============================================================
# Create TC06_SignUp as a copy of TC05_PlacingOrder (keeps all
# styles, merges, data validation and the hyperlink intact).
# ============================================================
$srcSheet = $wb.Worksheets.Item("TC05_PlacingOrder")
$srcSheet.Copy($null, $srcSheet)
$signUp = $wb.Worksheets.Item($wb.Worksheets.Count)
$signUp.Name = "TC06_SignUp"

# ============================================================
# Create TC07_Login as a copy of TC06_SignUp.
# ============================================================
$signUp.Copy($null, $signUp)
$login = $wb.Worksheets.Item($wb.Worksheets.Count)
$login.Name = "TC07_Login"

# ============================================================
# TC06_SignUp content
# ============================================================
$signUp.Range("D1").Value = "Sign up creating a user account"
$signUp.Range("B11").Value = "Click ""Sign Up"" button at navigational panel"
$signUp.Range("D11").Value = "Opens a ""Sign Up"" window"
$signUp.Range("B12").Value = "Click ""Sign Up"" button at the bottom of the window"
$signUp.Range("D12").Value = "Shows a subwindow with text ""Please fill out Username and Password"""
$signUp.Range("B13").Value = "Fill Username textbox and press ""Sign Up"" button"
$signUp.Range("D13").Value = "Shows a subwindow with text ""Please fill out Username and Password"""
$signUp.Range("B14").Value = "Empty Username textbox and fill Password textbox and press ""Sign Up"" button"
$signUp.Range("D14").Value = "Shows a subwindow with text ""Please fill out Username and Password"""
$signUp.Range("B15").Value = "Fill Username textbox that is already registered and press ""Sign Up"" button"
$signUp.Range("D15").Value = "Shows a subwindow ""This user already exist"""
$signUp.Range("B16").Value = "Fill Username textbox that is not registered and press ""Sign Up"" button"
$signUp.Range("D16").Value = "Shows a subwindow ""Sign up succsesfull"""
$signUp.Range("B17").Value = ""
$signUp.Range("D17").Value = ""

# Row heights that differ from the TC05_PlacingOrder template
$signUp.Rows.Item(11).AutoFit()
$signUp.Rows.Item(12).RowHeight = 44.25
$signUp.Rows.Item(13).RowHeight = 44.25
$signUp.Rows.Item(14).RowHeight = 51.75
$signUp.Rows.Item(15).RowHeight = 48.75
$signUp.Rows.Item(16).RowHeight = 45
$signUp.Rows.Item(17).RowHeight = 15.75

# ============================================================
# TC07_Login content
# ============================================================
$login.Range("B11").Value = "Click ""Log In"" button at navigational panel"
$login.Range("D11").Value = "Opens a ""Log In"" window"
$login.Range("B12").Value = "Click ""Log In"" button at the bottom of the window"
$login.Range("D12").Value = "Shows a subwindow with text ""Please fill out Username and Password"""
$login.Range("B13").Value = "Fill Username textbox and press ""Log In"" button"
$login.Range("D13").Value = "Shows a subwindow with text ""Please fill out Username and Password"""
$login.Range("B14").Value = "Empty Username textbox and fill Password textbox and press ""Log In"" button"
$login.Range("D14").Value = "Shows a subwindow with text ""Please fill out Username and Password"""
$login.Range("B15").Value = "Fill Username textbox that is already registered with incorrect password press ""Log In"" button"
$login.Range("D15").Value = "Shows a subwindow ""Wrong password"""
$login.Range("B16").Value = "Fill Username textbox that is already registered with correct password, press ""Log In"" button"
$login.Range("D16").Value = "Login succesfull. Page updates and Navigation Panel shows ""Welcome {username}"" and ""Log out"" button"
$login.Range("B17").Value = "Click ""Log Out"" button at navigational panel"
$login.Range("D17").Value = "Log out succesfull and website returns to the default state"
$login.Range("D1").Value = ""

# Row heights that differ from the TC05_PlacingOrder template
$login.Rows.Item(11).RowHeight = 36
$login.Rows.Item(12).RowHeight = 49.5
$login.Rows.Item(13).RowHeight = 54
$login.Rows.Item(14).RowHeight = 45
$login.Rows.Item(15).RowHeight = 50.25
$login.Rows.Item(16).RowHeight = 66
$login.Rows.Item(17).RowHeight = 45

# ============================================================
# Selections / active sheet (matches the final workbook state)
# ============================================================
$srcSheet.Activate()
$srcSheet.Range("L15").Select()

$signUp.Activate()
$signUp.Range("B11:E16").Select()

$login.Activate()
$login.Range("K16").Select()

